# Auto-generated script updating cached market-price columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled
# price-refresh run described in the commit message.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 117
$ws.Range("I33").Value = 122.916664
$ws.Range("K33").Value = 122.916664
$ws.Range("M33").Value = 106.083336
$ws.Range("H43").Value = 2572990.8
$ws.Range("I43").Value = 3086700
$ws.Range("K43").Value = 3086700
$ws.Range("M43").Value = -3086631
$ws.Range("H107").Value = 651.875
$ws.Range("I107").Value = 479.15384
$ws.Range("J107").Value = 1400.3334
$ws.Range("K107").Value = 479.15384
$ws.Range("L107").Value = 1400.3334
$ws.Range("M107").Value = 1440.84616
$ws.Range("N107").Value = -5240.3334
$ws.Range("H113").Value = 35851.227
$ws.Range("I113").Value = 65655.31
$ws.Range("J113").Value = 4060.2
$ws.Range("K113").Value = 65655.31
$ws.Range("L113").Value = 4060.2
$ws.Range("M113").Value = -62401.31
$ws.Range("N113").Value = -10568.2
$ws.Range("H135").Value = 3142.8647
$ws.Range("I135").Value = 1897.2333
$ws.Range("K135").Value = 17075.0997
$ws.Range("M135").Value = -14540.0997
$ws.Range("H138").Value = 6125.5415
$ws.Range("I138").Value = 2047.5
$ws.Range("J138").Value = 6496.273
$ws.Range("K138").Value = 6142.5
$ws.Range("L138").Value = 19488.819
$ws.Range("M138").Value = -1002.5
$ws.Range("N138").Value = -29768.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 543601.9399999999
$ws.Range("I2").Value = 738429.9
$ws.Range("J2").Value = 2413.111
$ws.Range("K2").Value = 738429.9
$ws.Range("L2").Value = 2413.111
$ws.Range("M2").Value = -738316.9
$ws.Range("N2").Value = -2639.111
$ws.Range("H61").Value = 21766.217
$ws.Range("I61").Value = 16491.3
$ws.Range("K61").Value = 16491.3
$ws.Range("M61").Value = -16279.3
$ws.Range("H74").Value = 6413903
$ws.Range("I74").Value = 11365073
$ws.Range("K74").Value = 11365073
$ws.Range("M74").Value = -11364199
$ws.Range("H77").Value = 6413903
$ws.Range("I77").Value = 11365073
$ws.Range("K77").Value = 56825365
$ws.Range("M77").Value = -56820997
$ws.Range("H116").Value = 543601.9399999999
$ws.Range("I116").Value = 738429.9
$ws.Range("J116").Value = 2413.111
$ws.Range("K116").Value = 738429.9
$ws.Range("L116").Value = 2413.111
$ws.Range("M116").Value = -736135.9
$ws.Range("N116").Value = -7001.111
$ws.Range("H122").Value = 429192
$ws.Range("I122").Value = 790213.9
$ws.Range("J122").Value = 7999.8335
$ws.Range("K122").Value = 2370641.7
$ws.Range("L122").Value = 23999.5005
$ws.Range("M122").Value = -2368191.7
$ws.Range("N122").Value = -28899.5005
$ws.Range("H132").Value = 35886.65
$ws.Range("I132").Value = 41250.2
$ws.Range("K132").Value = 123750.6
$ws.Range("M132").Value = -121220.6
$ws.Range("H136").Value = 21766.217
$ws.Range("I136").Value = 16491.3
$ws.Range("K136").Value = 49473.89999999999
$ws.Range("M136").Value = -46923.89999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 543601.9399999999
$ws.Range("I3").Value = 738429.9
$ws.Range("J3").Value = 2413.111
$ws.Range("K3").Value = 738429.9
$ws.Range("L3").Value = 2413.111
$ws.Range("M3").Value = -738315.9
$ws.Range("N3").Value = -2641.111
$ws.Range("H94").Value = 1337258.1
$ws.Range("I94").Value = 2100991.2
$ws.Range("K94").Value = 2100991.2
$ws.Range("M94").Value = -2100540.2
$ws.Range("H123").Value = 58799.332
$ws.Range("J123").Value = 58799.332
$ws.Range("L123").Value = 58799.332
$ws.Range("N123").Value = -68599.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1668001.5
$ws.Range("I58").Value = 1668001.5
$ws.Range("K58").Value = 1668001.5
$ws.Range("M58").Value = -1667798.5
$ws.Range("H132").Value = 17547810
$ws.Range("J132").Value = 27003.25
$ws.Range("L132").Value = 81009.75
$ws.Range("N132").Value = -86069.75
$ws.Range("H136").Value = 1668001.5
$ws.Range("I136").Value = 1668001.5
$ws.Range("K136").Value = 5004004.5
$ws.Range("M136").Value = -5001454.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 5005893
$ws.Range("I64").Value = 1999.5
$ws.Range("K64").Value = 5998.5
$ws.Range("M64").Value = -5728.5
$ws.Range("H67").Value = 5005893
$ws.Range("I67").Value = 1999.5
$ws.Range("K67").Value = 5998.5
$ws.Range("M67").Value = -5062.5
$ws.Range("H68").Value = 374440.75
$ws.Range("J68").Value = 560826.8
$ws.Range("L68").Value = 1682480.4
$ws.Range("N68").Value = -1684102.4
$ws.Range("H71").Value = 374440.75
$ws.Range("J71").Value = 560826.8
$ws.Range("L71").Value = 5047441.2
$ws.Range("N71").Value = -5055553.2
$ws.Range("H114").Value = 799.5294
$ws.Range("I114").Value = 808.3333
$ws.Range("K114").Value = 2424.9999
$ws.Range("M114").Value = 829.0001000000002
$ws.Range("H122").Value = 378.1875
$ws.Range("J122").Value = 485.33334
$ws.Range("L122").Value = 4368.00006
$ws.Range("N122").Value = -9268.00006
$ws.Range("H137").Value = 28672082
$ws.Range("I137").Value = 35715964
$ws.Range("J137").Value = 12236354
$ws.Range("K137").Value = 107147892
$ws.Range("L137").Value = 36709062
$ws.Range("M137").Value = -107142792
$ws.Range("N137").Value = -36719262
$ws.Range("H140").Value = 1343.4445
$ws.Range("I140").Value = 1343.4445
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 4030.3335
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 1149.6665
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 15928809
$ws.Range("I113").Value = 18583250
$ws.Range("J113").Value = 2166
$ws.Range("K113").Value = 18583250
$ws.Range("L113").Value = 2166
$ws.Range("M113").Value = -18581080
$ws.Range("N113").Value = -6506
$ws.Range("H122").Value = 529628.5600000001
$ws.Range("I122").Value = 790228.2
$ws.Range("J122").Value = 8429.286
$ws.Range("K122").Value = 2370684.6
$ws.Range("L122").Value = 25287.858
$ws.Range("M122").Value = -2368234.6
$ws.Range("N122").Value = -30187.858
$ws.Range("H126").Value = 3658
$ws.Range("I126").Value = 2111.3635
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 6334.0905
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -3864.0905
$ws.Range("N126").Value = -49940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4219.533
$ws.Range("I7").Value = 3729.372
$ws.Range("J7").Value = 5459.353
$ws.Range("K7").Value = 3729.372
$ws.Range("L7").Value = 5459.353
$ws.Range("M7").Value = -3617.372
$ws.Range("N7").Value = -5683.353
$ws.Range("H16").Value = 50053044
$ws.Range("I16").Value = 66734604
$ws.Range("K16").Value = 66734604
$ws.Range("M16").Value = -66734434
$ws.Range("H22").Value = 1347.4
$ws.Range("J22").Value = 803.5714
$ws.Range("L22").Value = 803.5714
$ws.Range("N22").Value = -1393.5714
$ws.Range("H27").Value = 1347.4
$ws.Range("J27").Value = 803.5714
$ws.Range("L27").Value = 803.5714
$ws.Range("N27").Value = -1017.5714
$ws.Range("H40").Value = 55560308
$ws.Range("J40").Value = 66673068
$ws.Range("L40").Value = 66673068
$ws.Range("N40").Value = -66673340
$ws.Range("H68").Value = 1516721.1
$ws.Range("I68").Value = 3248347.2
$ws.Range("J68").Value = 1548.125
$ws.Range("K68").Value = 3248347.2
$ws.Range("L68").Value = 1548.125
$ws.Range("M68").Value = -3247598.2
$ws.Range("N68").Value = -3046.125
$ws.Range("H71").Value = 1516721.1
$ws.Range("I71").Value = 3248347.2
$ws.Range("J71").Value = 1548.125
$ws.Range("K71").Value = 16241736
$ws.Range("L71").Value = 7740.625
$ws.Range("M71").Value = -16237992
$ws.Range("N71").Value = -15228.625
$ws.Range("H126").Value = 4219.533
$ws.Range("I126").Value = 3729.372
$ws.Range("J126").Value = 5459.353
$ws.Range("K126").Value = 11188.116
$ws.Range("L126").Value = 16378.059
$ws.Range("M126").Value = -8718.116
$ws.Range("N126").Value = -21318.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 770.5333000000001
$ws.Range("I113").Value = 662.6
$ws.Range("J113").Value = 986.4
$ws.Range("K113").Value = 1987.8
$ws.Range("L113").Value = 2959.2
$ws.Range("M113").Value = 182.1999999999998
$ws.Range("N113").Value = -7299.2
$ws.Range("H132").Value = 2323191.8
$ws.Range("I132").Value = 3482163.5
$ws.Range("K132").Value = 10446490.5
$ws.Range("M132").Value = -10443960.5

